$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.955.26"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = "'2.657.65"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'537.98"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.15%  '
$ws.Range('D6').Value = "'144.60"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.99%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = "'6.59"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').Value = "'0.337"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').Value = "'3.111.42"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('D14').Value = "'59.883.55"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('D15').Value = "'20.95"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').Value = "'2.641.29"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').Value = "'0.0000134"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = "'343.90"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('E19').Value = '  +2.09%  '
$ws.Range('D20').Value = "'10.22"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').Value = "'6.42"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = "'67.31"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').Value = "'0.413"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = "'0.997"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +3.05%  '
$ws.Range('D28').Value = "'0.0₃0752"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.35%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = "'1.66"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('D31').Value = "'5.87"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').Value = "'18.94"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('D34').Value = "'4.02"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').Value = "'0.840"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').Value = "'0.832"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.86%  '
$ws.Range('D39').Value = "'291.09"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.94%  '
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('D43').Value = "'10.75"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = "'0.0952"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = "'0.0533"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('D46').Value = "'1.974.96"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('D47').Value = "'18.66"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.59%  '
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('D50').Value = "'111.60"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('E51').Value = '  +0.19%  '
